$d = $word.ActiveDocument

# Locate the (currently empty) paragraph that precedes the existing
# red "Lorem ipsum..." block -- this is where the new blue paragraph
# block is inserted, per the diff (w14:paraId="52E27C42").
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "`r" -or $p.Range.Text -eq "") {
        $next = $null
        if ($i -lt $d.Paragraphs.Count) {
            $next = $d.Paragraphs.Item($i + 1)
        }
        if ($next -ne $null -and $next.Range.Text.StartsWith("Lorem ipsum dolor sit amet")) {
            $target = $p
            break
        }
    }
}

$w = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"
$w14 = "http://schemas.microsoft.com/office/word/2010/wordml"

$xml = '<w:p xmlns:w="' + $w + '" xmlns:w14="' + $w14 + '" w14:paraId="52E27C42" w14:textId="55C5DEB8" w:rsidR="00772C6C" w:rsidRDefault="00772C6C">' +
'<w:pPr><w:rPr><w:noProof/><w:color w:val="4472C4" w:themeColor="accent1"/><w:lang w:val="es-ES"/></w:rPr></w:pPr>' +
'<w:r><w:rPr><w:noProof/><w:color w:val="4472C4" w:themeColor="accent1"/><w:lang w:val="es-ES"/></w:rPr>' +
'<w:t>Lorem ipsum dolor sit amet, consectetuer adipiscing elit. Maecenas porttitor congue massa. Fusce posuere, magna sed pulvinar ultricies, purus lectus malesuada libero, sit amet commodo magna eros quis urna. Nunc viverra imperdiet enim. Fusce est.</w:t>' +
'</w:r></w:p>' +
'<w:p xmlns:w="' + $w + '"><w:pPr><w:rPr><w:noProof/><w:color w:val="4472C4" w:themeColor="accent1"/><w:lang w:val="es-ES"/></w:rPr></w:pPr>' +
'<w:r><w:rPr><w:noProof/><w:color w:val="4472C4" w:themeColor="accent1"/><w:lang w:val="es-ES"/></w:rPr>' +
'<w:t>Vivamus a tellus. Pellentesque habitant morbi tristique senectus et netus et malesuada fames ac turpis egestas. Proin pharetra nonummy pede. Mauris et orci. Aenean nec lorem.</w:t>' +
'</w:r></w:p>' +
'<w:p xmlns:w="' + $w + '"><w:pPr><w:rPr><w:noProof/><w:color w:val="4472C4" w:themeColor="accent1"/><w:lang w:val="es-ES"/></w:rPr></w:pPr>' +
'<w:r><w:rPr><w:noProof/><w:color w:val="4472C4" w:themeColor="accent1"/><w:lang w:val="es-ES"/></w:rPr>' +
'<w:t>In porttitor. Donec laoreet nonummy augue. Suspendisse dui purus, scelerisque at, vulputate vitae, pretium mattis, nunc. Mauris eget neque at sem venenatis eleifend. Ut nonummy.</w:t>' +
'</w:r></w:p>' +
'<w:p xmlns:w="' + $w + '"><w:pPr><w:rPr><w:noProof/><w:color w:val="4472C4" w:themeColor="accent1"/><w:lang w:val="es-ES"/></w:rPr></w:pPr>' +
'<w:r><w:rPr><w:noProof/><w:color w:val="4472C4" w:themeColor="accent1"/><w:lang w:val="es-ES"/></w:rPr>' +
'<w:t>Fusce aliquet pede non pede. Suspendisse dapibus lorem pellentesque magna. Integer nulla. Donec blandit feugiat ligula. Donec hendrerit, felis et imperdiet euismod, purus ipsum pretium metus, in lacinia nulla nisl eget sapien.</w:t>' +
'</w:r></w:p>'

$target.Range.InsertXML($xml)

Write-Output "inserted blue paragraph block"
